$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.778.13'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '2.497.89'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("E9").Value = '  +5.10%  '
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("E11").Value = '  +3.98%  '
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("D14").Value = '2.913.63'
$ws.Range("D15").Value = '67.642.07'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '2.444.65'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '352.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.07%  '
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").Value = '2.624.11'
$ws.Range("E27").Value = '  +1.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = '0.0₃0913'
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '511.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.49%  '
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +7.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("E40").Value = '  +6.19%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.331'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("E44").Value = '  +3.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.15%  '
$ws.Range("E47").Value = '  +4.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.516'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0746'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.587'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.32%  '
